# Refresh the cryptocurrency price/volume snapshot (GitHub Actions bot run).
# Updates Price (D) and Volume(1h) (E) figures for most rows, and for rows
# 40/41 also swaps which coin (RenderToken / InjectiveProtocol) occupies
# which row, along with its Link/Price/Volume.
#
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as literal text (e.g. "34.509.89", "1.00", "0.604")
# instead of silently reinterpreting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'34.509.89"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3
$ws.Range("D3").Value = "'1.802.27"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'224.48"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6
$ws.Range("D6").Value = "'0.604"
$ws.Range("E6").Value = "  +0.94%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").Value = "'42.24"
$ws.Range("E8").Value = "  +17.26%  "

# Row 9
$ws.Range("D9").Value = "'0.291"
$ws.Range("E9").Value = "  +0.37%  "

# Row 10
$ws.Range("E10").Value = "  -1.79%  "

# Row 11
$ws.Range("D11").Value = "'0.0994"
$ws.Range("E11").Value = "  +2.96%  "

# Row 12
$ws.Range("D12").Value = "'2.063.90"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").Value = "'1.795.28"
$ws.Range("E13").Value = "  -0.47%  "

# Row 14
$ws.Range("E14").Value = "  -2.58%  "

# Row 15
$ws.Range("D15").Value = "'34.479.66"
$ws.Range("E15").Value = "  +0.32%  "

# Row 16
$ws.Range("D16").Value = "'0.627"
$ws.Range("E16").Value = "  -0.27%  "

# Row 17
$ws.Range("E17").Value = "  -0.57%  "

# Row 18
$ws.Range("D18").Value = "'67.26"
$ws.Range("E18").Value = "  -1.92%  "

# Row 19
$ws.Range("D19").Value = "'240.30"
$ws.Range("E19").Value = "  -0.88%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0767"
$ws.Range("E20").Value = "  -0.92%  "

# Row 21
$ws.Range("E21").Value = "  -1.36%  "

# Row 22
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("E23").Value = "  +6.30%  "

# Row 24
$ws.Range("E24").Value = "  -3.04%  "

# Row 25
$ws.Range("D25").Value = "'171.97"
$ws.Range("E25").Value = "  +0.51%  "

# Row 26
$ws.Range("D26").Value = "'7.65"
$ws.Range("E26").Value = "  -2.42%  "

# Row 27
$ws.Range("E27").Value = "  +0.19%  "

# Row 28
$ws.Range("E28").Value = "  +0.12%  "

# Row 30
$ws.Range("E30").Value = "  +0.27%  "

# Row 31
$ws.Range("D31").Value = "'1.22"
$ws.Range("E31").Value = "  -0.39%  "

# Row 32
$ws.Range("D32").Value = "'3.87"
$ws.Range("E32").Value = "  -0.96%  "

# Row 33
$ws.Range("D33").Value = "'0.0513"
$ws.Range("E33").Value = "  -0.43%  "

# Row 34
$ws.Range("D34").Value = "'1.78"
$ws.Range("E34").Value = "  +0.92%  "

# Row 35
$ws.Range("D35").Value = "'87.28"
$ws.Range("E35").Value = "  +7.55%  "

# Row 36
$ws.Range("E36").Value = "  -0.35%  "

# Row 37
$ws.Range("D37").Value = "'1.316.80"
$ws.Range("E37").Value = "  -3.44%  "

# Row 38
$ws.Range("E38").Value = "  -0.50%  "

# Row 39
$ws.Range("E39").Value = "  +0.58%  "

# Row 40
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'14.71"
$ws.Range("E40").Value = "  +11.24%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.33"
$ws.Range("E41").Value = "  -1.07%  "

# Row 42
$ws.Range("E42").Value = "  +5.36%  "

# Row 43
$ws.Range("D43").Value = "'2.42"
$ws.Range("E43").Value = "  +0.28%  "

# Row 45
$ws.Range("D45").Value = "'0.937"
$ws.Range("E45").Value = "  +0.13%  "

# Row 46
$ws.Range("D46").Value = "'0.0519"
$ws.Range("E46").Value = "  +3.51%  "

# Row 47
$ws.Range("D47").Value = "'1.965.42"
$ws.Range("E47").Value = "  +0.09%  "

# Row 48
$ws.Range("E48").Value = "  -0.01%  "

# Row 49
$ws.Range("E49").Value = "  +0.11%  "

# Row 50
$ws.Range("D50").Value = "'100.49"
$ws.Range("E50").Value = "  -1.23%  "

# Row 51
$ws.Range("D51").Value = "'0.0607"
$ws.Range("E51").Value = "  +0.46%  "
